# Apply "Update with Correct Forecast output" edit to the workbook.
#
# Summary of changes on the "Forecast Comparison" sheet:
#   - Insert a new column B "Week_Start_Date" (everything from old column B
#     onward shifts one column to the right).
#   - Column A ("Week") values lose their zero padding: W01 -> W1 ... W16.
#   - New column B gets the ISO week-start date for each row (stored as text).
#   - Column D ("MyForecast", old column C) gets corrected values for a few
#     rows.
#   - Column J ("is_holiday_week", old column I) switches from numeric 0/1 to
#     boolean FALSE/TRUE.
#
# The "Summary" sheet has a handful of values that need to be refreshed to
# stay consistent with the corrected forecast numbers above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1) Insert the new "Week_Start_Date" column before the old column B (ASIN).
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# 2) Row data: Week, Week_Start_Date, MyForecast (col D)
$weeks = @(
    @{ Row = 2;  Week = "W1";  Date = "2025-01-05"; MyForecast = 17 },
    @{ Row = 3;  Week = "W2";  Date = "2025-01-12"; MyForecast = 19 },
    @{ Row = 4;  Week = "W3";  Date = "2025-01-19"; MyForecast = 17 },
    @{ Row = 5;  Week = "W4";  Date = "2025-01-26"; MyForecast = 19 },
    @{ Row = 6;  Week = "W5";  Date = "2025-02-02"; MyForecast = 20 },
    @{ Row = 7;  Week = "W6";  Date = "2025-02-09"; MyForecast = 19 },
    @{ Row = 8;  Week = "W7";  Date = "2025-02-16"; MyForecast = 22 },
    @{ Row = 9;  Week = "W8";  Date = "2025-02-23"; MyForecast = 23 },
    @{ Row = 10; Week = "W9";  Date = "2025-03-02"; MyForecast = 17 },
    @{ Row = 11; Week = "W10"; Date = "2025-03-09"; MyForecast = 23 },
    @{ Row = 12; Week = "W11"; Date = "2025-03-16"; MyForecast = 23 },
    @{ Row = 13; Week = "W12"; Date = "2025-03-23"; MyForecast = 20 },
    @{ Row = 14; Week = "W13"; Date = "2025-03-30"; MyForecast = 20 },
    @{ Row = 15; Week = "W14"; Date = "2025-04-06"; MyForecast = 18 },
    @{ Row = 16; Week = "W15"; Date = "2025-04-13"; MyForecast = 17 },
    @{ Row = 17; Week = "W16"; Date = "2025-04-20"; MyForecast = 20 }
)

foreach ($w in $weeks) {
    $r = $w.Row
    $ws.Cells.Item($r, 1).Value = $w.Week
    # Leading apostrophe forces the date-looking text to stay a plain string
    # instead of being auto-converted into a date serial number.
    $ws.Cells.Item($r, 2).Value = "'" + $w.Date
    $ws.Cells.Item($r, 4).Value = $w.MyForecast
    # Column J (is_holiday_week) becomes a real boolean.
    $ws.Cells.Item($r, 10).Value = $false
}

# 3) Refresh the dependent metrics on the "Summary" sheet.
# All Value column entries on this sheet are plain text, so force the
# look-like-a-number / look-like-a-date values to stay text with a leading
# apostrophe (otherwise Excel would silently convert them to a number or a
# date serial).
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(9, 2).Value = "'314"
$summary.Cells.Item(10, 2).Value = "'156"
$summary.Cells.Item(11, 2).Value = "'72"
$summary.Cells.Item(13, 2).Value = "'2025-02-23"
